# Apply "excel report vivado aggiornato" update to the "Data" sheet.
# Re-run values for solutions #5 (row 9), #6 (row 10), #7 (row 12) and
# #8 (row 13). Columns K (sum), M and Z are driven by existing formulas
# in the sheet, so they recompute automatically once D:J / O:Y are set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---- Row 9 (Solution "5") --------------------------------------------
$ws.Range("D9").Formula = "=1000*0.00124085135757923"
$ws.Range("E9").Formula = "=1000*0.0000825246315798722"
$ws.Range("F9").Formula = "=1000*0.000960682868026197"
$ws.Range("G9").Formula = "=1000*0.000272355420747772"
$ws.Range("H9").Formula = "=1000*0.000266662013018504"
$ws.Range("I9").Formula = "=1000*0.00000428146995545831"
$ws.Range("J9").Formula = "=1000*0.000425589532824233"

$ws.Range("O9").Value = 98
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 72
$ws.Range("R9").Value = 1
$ws.Range("W9").Value = 57
$ws.Range("Y9").Value = 4.33

# ---- Row 10 (Solution "6") ---------------------------------------------
$ws.Range("E10").Formula = "=1000*0.000325270200846717"
$ws.Range("F10").Formula = "=1000*0.00235283561050892"
$ws.Range("G10").Formula = "=1000*0.000263715046457946"
$ws.Range("H10").Formula = "=1000*0.000575191108509898"
$ws.Range("I10").Formula = "=1000*0.00000701051294527133"
$ws.Range("J10").Formula = "=1000*0.000750690058339387"

$ws.Range("O10").Value = 413
$ws.Range("Q10").Value = 843
$ws.Range("S10").Value = 2
$ws.Range("W10").Value = 66
$ws.Range("Y10").Value = 3.469

# ---- Row 12 (Solution "7") ----------------------------------------------
$ws.Range("D12").Formula = "=1000*0.00123103871010244"
$ws.Range("E12").Formula = "=1000*0.000106978332041763"
$ws.Range("F12").Formula = "=1000*0.000972227076999843"
$ws.Range("G12").Formula = "=1000*0.000247065967414528"
$ws.Range("H12").Formula = "=1000*0.000258892890997231"
$ws.Range("I12").Formula = "=1000*0.00000258562499766413"
$ws.Range("J12").Formula = "=1000*0.00041881192009896"

$ws.Range("O12").Value = 145
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 93
$ws.Range("R12").Value = 1
$ws.Range("W12").Value = 61
$ws.Range("Y12").Value = 4.33

# ---- Row 13 (Solution "8") ----------------------------------------------
$ws.Range("E13").Formula = "=1000*0.000328624591929838"
$ws.Range("F13").Formula = "=1000*0.00256039062514901"
$ws.Range("G13").Formula = "=1000*0.00029804851510562"
$ws.Range("H13").Formula = "=1000*0.00114636321086437"
$ws.Range("I13").Formula = "=1000*0.00000300649980999879"
$ws.Range("J13").Formula = "=1000*0.00132495700381696"

$ws.Range("O13").Value = 1145
$ws.Range("Q13").Value = 864
$ws.Range("S13").Value = 2
$ws.Range("W13").Value = 59
$ws.Range("Y13").Value = 3.097

# ---- Selection moved to Z14 on the "Data" sheet -------------------------
[void]$ws.Range("Z14").Select()
